# Incident_Management_TestData.xlsx edit
# - new Incident IDs (row 2 & row 3) reflecting latest test run
# - brand new row 4 dataset for a new test (testResolveIncidentTicket)
# - four new trailing columns: Configuration Item / Cause Code / Sub Cause Code /
#   Mitigation and Solution Steps (+ one new blank trailing column to keep the
#   right-hand border formatting consistent)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Write all new / changed cell VALUES first (formats are copied afterwards
#    so a PasteSpecial(formats) never clobbers a value we just set).
# ---------------------------------------------------------------------------

# Header row - four new columns
$ws.Range("N1").Value = "Configuration Item"
$ws.Range("O1").Value = "Cause Code"
$ws.Range("P1").Value = "Sub Cause Code"
$ws.Range("Q1").Value = "Mitigation and Solution Steps"

# Row 2 - refreshed incident id
$ws.Range("C2").Value = "INC0021040"

# Row 3 - refreshed incident id
$ws.Range("C3").Value = "INC0020960"

# Row 4 - brand new data row for testResolveIncidentTicket
$ws.Range("B4").Value = "testResolveIncidentTicket"
$ws.Range("C4").Value = "INC0021040"
$ws.Range("E4").Value = "Passed"
$ws.Range("H4").Value = "Sathyanarayanan V"
# Leading apostrophe -> Excel "quote prefix" (text-stored-as-text) cell, matching
# the authored workbook's quotePrefix="1" style for this cell.
$ws.Range("N4").Value = "'162LOAN-TPL-G"
$ws.Range("O4").Value = "Change"
$ws.Range("P4").Value = "Change Failure"
$ws.Range("Q4").Value = "Mitigation and Solution Steps for "

# ---------------------------------------------------------------------------
# 2) Copy cell formatting (without values) from stable, unmodified cells that
#    already carry the style each destination cell needs to end up with.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $src = $ws.Range($srcAddr)
    $src.Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# Row 1 new header cells
Copy-Format "O1" "P1"   # P1 must take O1's CURRENT style (8) before O1 itself is restyled
Copy-Format "A1" "O1"   # O1 -> style 1
Copy-Format "M1" "N1"   # N1 -> style 2
Copy-Format "R1" "S1"   # S1 (new trailing column) -> style 7

# Row 2
Copy-Format "B2" "P2"   # P2 -> style 5
Copy-Format "Q2" "S2"   # S2 (new trailing column) -> style 3

# Row 3
Copy-Format "B2" "P3"   # P3 -> style 5
Copy-Format "Q3" "S3"   # S3 (new trailing column) -> style 3

# Row 4
Copy-Format "C2" "C4"   # C4 -> default style (no explicit style id)
Copy-Format "E2" "E4"   # E4 -> default style (no explicit style id)
Copy-Format "B2" "P4"   # P4 -> style 5
Copy-Format "R4" "S4"   # S4 (new trailing column) -> style 3

# Row 5
Copy-Format "B2" "P5"   # P5 -> style 5
Copy-Format "Q5" "S5"   # S5 (new trailing column) -> style 3

# Row 6
Copy-Format "B2" "P6"   # P6 -> style 5
Copy-Format "Q6" "S6"   # S6 (new trailing column) -> style 3

# ---------------------------------------------------------------------------
# 3) New trailing-column widths (bestFit-style widths recomputed by Excel for
#    the newly added / shifted columns N..S).
# ---------------------------------------------------------------------------

$ws.Columns.Item(14).ColumnWidth = 16.92   # N - Configuration Item
$ws.Columns.Item(15).ColumnWidth = 9.59    # O - Cause Code
$ws.Columns.Item(16).ColumnWidth = 13.25   # P - Sub Cause Code
$ws.Columns.Item(17).ColumnWidth = 27.59   # Q - Mitigation and Solution Steps
$ws.Columns.Item(18).ColumnWidth = 17.09   # R
$ws.Columns.Item(19).ColumnWidth = 16.42   # S

# ---------------------------------------------------------------------------
# 4) Selection, matching the authored workbook's saved cursor position.
# ---------------------------------------------------------------------------

$ws.Range("F11").Select()
